# Daily refresh of the BP "terminal gate pricing" workbook.
#
# Each state/product section keeps a rolling 2-day window: the newest
# "Effective Date" row(s) become the prior day's row(s), and a brand new
# "today" row is written in with its own prices. Columns are:
#   A = Effective Date (serial date number)
#   B = (blank)
#   C = Terminal (unchanged - shared string, not touched here)
#   D = Diesel, E = ULP, F = PULP, G = e10
#
# Below, every data row in the sheet is updated in place: column A gets its
# new date serial, and D/E/F/G get their new prices (only the columns that
# actually have a value in that row are touched; "N/A" text cells and blank
# e10 cells are left alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ col = value; ... }
$updates = [ordered]@{
    8  = @{ A = 45986; D = 173.02; E = 161.29; F = 171.29; G = 161.46 }
    9  = @{ A = 45986; D = 173.02; E = 161.29; F = 171.29; G = 161.46 }
    10 = @{ A = 45986; D = 174.97; E = 163.77; F = 173.77; G = 164.26 }
    11 = @{ A = 45983; D = 173.3;  E = 161.74; F = 171.74; G = 161.9  }
    12 = @{ A = 45983; D = 173.3;  E = 161.74; F = 171.74; G = 161.9  }
    13 = @{ A = 45983; D = 175.43; E = 164.39; F = 174.39; G = 164.88 }

    17 = @{ A = 45986; D = 178.25; E = 166.26; F = 176.26 }
    18 = @{ A = 45983; D = 178.72; E = 166.89; F = 176.89 }

    22 = @{ A = 45986; D = 174.14; E = 162.83; F = 172.43; G = 164.11 }
    23 = @{ A = 45986; D = 179.76; E = 167.59; F = 177.59 }
    24 = @{ A = 45986; D = 179.55; E = 167.97; F = 177.97 }
    25 = @{ A = 45986; D = 180.37; E = 167.38; F = 177.38; G = 167.42 }
    26 = @{ A = 45986; D = 179.06; E = 168.83; F = 178.83 }
    27 = @{ A = 45983; D = 174.43; E = 163.44; F = 173.04; G = 164.73 }
    28 = @{ A = 45983; D = 180.22; E = 168.2;  F = 178.2  }
    29 = @{ A = 45983; D = 180.01; E = 168.58; F = 178.58 }
    30 = @{ A = 45983; D = 180.83; E = 167.99; F = 177.99; G = 168.03 }
    31 = @{ A = 45983; D = 179.52; E = 169.45; F = 179.45 }

    35 = @{ A = 45986; D = 173.32; E = 160.76; F = 169.76 }
    36 = @{ A = 45983; D = 173.67; E = 161.38; F = 170.38 }

    40 = @{ A = 45986; D = 179;    E = 166.26; F = 176.26 }
    41 = @{ A = 45986; D = 178.7;  E = 166.68; F = 176.68 }
    42 = @{ A = 45983; D = 179.45; E = 166.87; F = 176.87 }
    43 = @{ A = 45983; D = 179.15; E = 167.29; F = 177.29 }

    47 = @{ A = 45986; D = 173.02; E = 162.59; F = 172.59 }
    48 = @{ A = 45986; D = 172.98; E = 162.74; F = 172.74 }
    49 = @{ A = 45983; D = 172.79; E = 162.81; F = 172.81 }
    50 = @{ A = 45983; D = 172.75; E = 162.96; F = 172.96 }

    54 = @{ A = 45986; D = 189.05; E = 177.14; F = 187.14 }
    55 = @{ A = 45986; D = 176.95; E = 173.76; F = 183.76 }
    56 = @{ A = 45986; D = 179.13 }
    57 = @{ A = 45986; D = 178.56; E = 168.03 }
    58 = @{ A = 45986; D = 174.47; E = 164.08; F = 174.08 }
    59 = @{ A = 45986; D = 181.2;  E = 175.06 }
    60 = @{ A = 45983; D = 189.5;  E = 177.72; F = 187.72 }
    61 = @{ A = 45983; D = 177.19; E = 174.39; F = 184.39 }
    62 = @{ A = 45983; D = 179.59 }
    63 = @{ A = 45983; D = 179.03; E = 168.65 }
    64 = @{ A = 45983; D = 174.94; E = 164.71; F = 174.71 }
    65 = @{ A = 45983; D = 181.68; E = 175.66 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
